$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wong3")

$ws.Cells.Item(2, 3).Value = 521
$ws.Cells.Item(2, 4).Value = 53.93374741200828
$ws.Cells.Item(3, 3).Value = 16
$ws.Cells.Item(3, 4).Value = 1.656314699792961
$ws.Cells.Item(4, 3).Value = 793
$ws.Cells.Item(4, 4).Value = 82.09109730848861
$ws.Cells.Item(5, 3).Value = 802
$ws.Cells.Item(5, 4).Value = 83.02277432712215
$ws.Cells.Item(6, 3).Value = 648
$ws.Cells.Item(6, 4).Value = 67.08074534161491
$ws.Cells.Item(7, 3).Value = 740
$ws.Cells.Item(7, 4).Value = 76.60455486542443
$ws.Cells.Item(8, 3).Value = 430
$ws.Cells.Item(8, 4).Value = 44.51345755693582
$ws.Cells.Item(9, 3).Value = 205
$ws.Cells.Item(9, 4).Value = 21.22153209109731
$ws.Cells.Item(10, 3).Value = 551
$ws.Cells.Item(10, 4).Value = 57.03933747412007
$ws.Cells.Item(11, 3).Value = 731
$ws.Cells.Item(11, 4).Value = 75.67287784679088
$ws.Cells.Item(12, 3).Value = 553
$ws.Cells.Item(12, 4).Value = 57.2463768115942
$ws.Cells.Item(13, 3).Value = 677
$ws.Cells.Item(13, 4).Value = 70.08281573498965
$ws.Cells.Item(14, 3).Value = 116
$ws.Cells.Item(14, 4).Value = 12.00828157349896
$ws.Cells.Item(15, 3).Value = 806
$ws.Cells.Item(15, 4).Value = 83.43685300207039
$ws.Cells.Item(16, 3).Value = 921
$ws.Cells.Item(16, 4).Value = 95.34161490683229
$ws.Cells.Item(18, 3).Value = 137
$ws.Cells.Item(18, 4).Value = 14.18219461697723
$ws.Cells.Item(19, 3).Value = 218
$ws.Cells.Item(19, 4).Value = 22.56728778467909
$ws.Cells.Item(20, 3).Value = 775
$ws.Cells.Item(20, 4).Value = 80.22774327122153
$ws.Cells.Item(21, 3).Value = 449
$ws.Cells.Item(21, 4).Value = 46.48033126293996
$ws.Cells.Item(22, 3).Value = 820
$ws.Cells.Item(22, 4).Value = 84.88612836438924
$ws.Cells.Item(23, 3).Value = 19
$ws.Cells.Item(23, 4).Value = 1.966873706004141
$ws.Cells.Item(24, 3).Value = 690
$ws.Cells.Item(24, 4).Value = 71.42857142857143
$ws.Cells.Item(26, 3).Value = 500
$ws.Cells.Item(26, 4).Value = 51.75983436853002
$ws.Cells.Item(28, 3).Value = 228
$ws.Cells.Item(28, 4).Value = 23.60248447204969
$ws.Cells.Item(29, 3).Value = 705
$ws.Cells.Item(29, 4).Value = 72.98136645962732
$ws.Cells.Item(30, 3).Value = 701
$ws.Cells.Item(30, 4).Value = 72.5672877846791
$ws.Cells.Item(31, 3).Value = 660
$ws.Cells.Item(31, 4).Value = 68.32298136645963
$ws.Cells.Item(32, 3).Value = 351
$ws.Cells.Item(32, 4).Value = 36.33540372670808
$ws.Cells.Item(33, 3).Value = 393
$ws.Cells.Item(33, 4).Value = 40.6832298136646
$ws.Cells.Item(34, 3).Value = 744
$ws.Cells.Item(34, 4).Value = 77.01863354037268
$ws.Cells.Item(35, 3).Value = 826
$ws.Cells.Item(35, 4).Value = 85.50724637681159
$ws.Cells.Item(36, 3).Value = 826
$ws.Cells.Item(36, 4).Value = 85.50724637681159
$ws.Cells.Item(37, 3).Value = 878
$ws.Cells.Item(37, 4).Value = 90.89026915113871
$ws.Cells.Item(38, 3).Value = 444
$ws.Cells.Item(38, 4).Value = 45.96273291925466
$ws.Cells.Item(39, 3).Value = 374
$ws.Cells.Item(39, 4).Value = 38.71635610766045
$ws.Cells.Item(41, 3).Value = 188
$ws.Cells.Item(41, 4).Value = 19.46169772256729
$ws.Cells.Item(42, 3).Value = 367
$ws.Cells.Item(42, 4).Value = 37.99171842650104
$ws.Cells.Item(43, 3).Value = 374
$ws.Cells.Item(43, 4).Value = 38.71635610766045
$ws.Cells.Item(44, 3).Value = 531
$ws.Cells.Item(44, 4).Value = 54.96894409937888
$ws.Cells.Item(45, 3).Value = 727
$ws.Cells.Item(45, 4).Value = 75.25879917184265
$ws.Cells.Item(46, 3).Value = 146
$ws.Cells.Item(46, 4).Value = 15.11387163561077
$ws.Cells.Item(47, 3).Value = 372
$ws.Cells.Item(47, 4).Value = 38.50931677018634
$ws.Cells.Item(48, 3).Value = 872
$ws.Cells.Item(48, 4).Value = 90.26915113871635
$ws.Cells.Item(49, 3).Value = 752
$ws.Cells.Item(49, 4).Value = 77.84679089026915
$ws.Cells.Item(50, 3).Value = 887
$ws.Cells.Item(50, 4).Value = 91.82194616977226
$ws.Cells.Item(51, 3).Value = 623
$ws.Cells.Item(51, 4).Value = 64.49275362318841
$ws.Cells.Item(52, 3).Value = 721
$ws.Cells.Item(52, 4).Value = 74.63768115942028
$ws.Cells.Item(53, 3).Value = 747
$ws.Cells.Item(53, 4).Value = 77.32919254658384
$ws.Cells.Item(54, 3).Value = 468
$ws.Cells.Item(54, 4).Value = 48.4472049689441
$ws.Cells.Item(55, 3).Value = 27
$ws.Cells.Item(55, 4).Value = 2.795031055900621
$ws.Cells.Item(56, 3).Value = 281
$ws.Cells.Item(56, 4).Value = 29.08902691511387
$ws.Cells.Item(57, 3).Value = 473
$ws.Cells.Item(57, 4).Value = 48.9648033126294
$ws.Cells.Item(58, 3).Value = 752
$ws.Cells.Item(58, 4).Value = 77.84679089026915
$ws.Cells.Item(59, 3).Value = 655
$ws.Cells.Item(59, 4).Value = 67.80538302277432
$ws.Cells.Item(60, 3).Value = 769
$ws.Cells.Item(60, 4).Value = 79.60662525879917
$ws.Cells.Item(61, 3).Value = 13
$ws.Cells.Item(61, 4).Value = 1.34575569358178
$ws.Cells.Item(62, 3).Value = 752
$ws.Cells.Item(62, 4).Value = 77.84679089026915
$ws.Cells.Item(64, 3).Value = 768
$ws.Cells.Item(64, 4).Value = 79.50310559006211
$ws.Cells.Item(65, 3).Value = 14
$ws.Cells.Item(65, 4).Value = 1.449275362318841
$ws.Cells.Item(66, 3).Value = 536
$ws.Cells.Item(66, 4).Value = 55.48654244306418
$ws.Cells.Item(67, 3).Value = 343
$ws.Cells.Item(67, 4).Value = 35.50724637681159
$ws.Cells.Item(68, 3).Value = 395
$ws.Cells.Item(68, 4).Value = 40.89026915113871
$ws.Cells.Item(69, 3).Value = 13
$ws.Cells.Item(69, 4).Value = 1.34575569358178
$ws.Cells.Item(70, 3).Value = 620
$ws.Cells.Item(70, 4).Value = 64.18219461697723
$ws.Cells.Item(71, 3).Value = 171
$ws.Cells.Item(71, 4).Value = 17.70186335403727
$ws.Cells.Item(72, 3).Value = 721
$ws.Cells.Item(72, 4).Value = 74.63768115942028
$ws.Cells.Item(73, 3).Value = 13
$ws.Cells.Item(73, 4).Value = 1.341589267285862
$ws.Cells.Item(74, 3).Value = 527
$ws.Cells.Item(74, 4).Value = 54.55486542443064
$ws.Cells.Item(75, 3).Value = 638
$ws.Cells.Item(75, 4).Value = 66.0455486542443
$ws.Cells.Item(76, 3).Value = 424
$ws.Cells.Item(76, 4).Value = 43.89233954451346
$ws.Cells.Item(77, 3).Value = 834
$ws.Cells.Item(77, 4).Value = 86.33540372670807
$ws.Cells.Item(78, 3).Value = 15
$ws.Cells.Item(78, 4).Value = 1.552795031055901
$ws.Cells.Item(79, 3).Value = 191
$ws.Cells.Item(79, 4).Value = 19.77225672877847
$ws.Cells.Item(80, 3).Value = 17
$ws.Cells.Item(80, 4).Value = 1.759834368530021
$ws.Cells.Item(81, 3).Value = 801
$ws.Cells.Item(81, 4).Value = 82.91925465838509
$ws.Cells.Item(82, 3).Value = 579
$ws.Cells.Item(82, 4).Value = 59.93788819875776
$ws.Cells.Item(84, 3).Value = 443
$ws.Cells.Item(84, 4).Value = 45.8592132505176
$ws.Cells.Item(85, 3).Value = 26
$ws.Cells.Item(85, 4).Value = 2.688728024819028
$ws.Cells.Item(86, 3).Value = 852
$ws.Cells.Item(86, 4).Value = 88.19875776397515
$ws.Cells.Item(87, 3).Value = 790
$ws.Cells.Item(87, 4).Value = 81.78053830227742
$ws.Cells.Item(88, 3).Value = 673
$ws.Cells.Item(88, 4).Value = 69.66873706004141
$ws.Cells.Item(89, 3).Value = 461
$ws.Cells.Item(89, 4).Value = 47.72256728778468
$ws.Cells.Item(90, 3).Value = 319
$ws.Cells.Item(90, 4).Value = 33.02277432712215
$ws.Cells.Item(91, 3).Value = 520
$ws.Cells.Item(91, 4).Value = 53.83022774327122
$ws.Cells.Item(92, 3).Value = 260
$ws.Cells.Item(92, 4).Value = 26.91511387163561
$ws.Cells.Item(93, 3).Value = 383
$ws.Cells.Item(93, 4).Value = 39.648033126294
$ws.Cells.Item(94, 3).Value = 198
$ws.Cells.Item(94, 4).Value = 20.47569803516029
$ws.Cells.Item(95, 3).Value = 13
$ws.Cells.Item(95, 4).Value = 1.34575569358178
$ws.Cells.Item(96, 3).Value = 542
$ws.Cells.Item(96, 4).Value = 56.10766045548654
$ws.Cells.Item(97, 3).Value = 796
$ws.Cells.Item(97, 4).Value = 82.40165631469979
$ws.Cells.Item(100, 3).Value = 167
$ws.Cells.Item(100, 4).Value = 17.28778467908903
$ws.Cells.Item(101, 3).Value = 4
$ws.Cells.Item(101, 4).Value = 0.4136504653567736
$ws.Cells.Item(102, 3).Value = 710
$ws.Cells.Item(102, 4).Value = 73.49896480331263
$ws.Cells.Item(105, 3).Value = 677
$ws.Cells.Item(105, 4).Value = 70.08281573498965
$ws.Cells.Item(106, 3).Value = 838
$ws.Cells.Item(106, 4).Value = 86.74948240165632

$ws.Name = "Euclid"

